$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C65").Value = "geo_mean_test"
$ws.Range("A65").Value = "Geo_Mean"
$ws.Range("B65").Value = "Test geometric mean"

$ws.Range("A67").Select()
